# Weekly crime-data refresh (NYPD CompStat Citywide report)
# - Bumps the "Volume/Number" and "Report Covering the Week" banner text
#   in-place (preserving the existing rich-text run formatting) by editing
#   only the substrings that changed.
# - Overwrites the weekly/28-day/YTD/2-year crime-complaint figures (rows
#   14-33) with the newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Banner text (A8): "Volume 31   Number  10" -> "...  11" ---
$a8 = $ws.Range("A8")
$a8text = $a8.Value2
$numStart = $a8text.LastIndexOf("10") + 1
$a8.Characters($numStart, 2).Text = "11"

# --- Banner text (C9): week-of dates ---
$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "3/11/2024"
$c9.Characters(47, 9).Text = "3/17/2024"

# --- Crime-complaint data table (rows 14-33) ---
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 28.571428571428
$ws.Range("F14").Value = 19
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = -24
$ws.Range("I14").Value = 65
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = -18.75
$ws.Range("L14").Value = -26.136363636363
$ws.Range("M14").Value = -29.347826086956
$ws.Range("N14").Value = -83.990147783251
$ws.Range("C15").Value = 41
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = 36.666666666666
$ws.Range("F15").Value = 129
$ws.Range("G15").Value = 100
$ws.Range("H15").Value = 29
$ws.Range("I15").Value = 303
$ws.Range("J15").Value = 315
$ws.Range("K15").Value = -3.809523809523
$ws.Range("L15").Value = -13.920454545454
$ws.Range("M15").Value = 28.936170212766
$ws.Range("N15").Value = -51.364365971107
$ws.Range("C16").Value = 278
$ws.Range("D16").Value = 256
$ws.Range("E16").Value = 8.59375
$ws.Range("F16").Value = 1157
$ws.Range("G16").Value = 1081
$ws.Range("H16").Value = 7.030527289546
$ws.Range("I16").Value = 3348
$ws.Range("J16").Value = 3163
$ws.Range("K16").Value = 5.848877647802
$ws.Range("L16").Value = 3.589108910891
$ws.Range("M16").Value = -8.047239769294
$ws.Range("N16").Value = -80.917640353377
$ws.Range("C17").Value = 563
$ws.Range("D17").Value = 471
$ws.Range("E17").Value = 19.532908704883
$ws.Range("F17").Value = 1941
$ws.Range("G17").Value = 1812
$ws.Range("H17").Value = 7.119205298013
$ws.Range("I17").Value = 5356
$ws.Range("J17").Value = 5158
$ws.Range("K17").Value = 3.838697169445
$ws.Range("L17").Value = 16.106655105137
$ws.Range("M17").Value = 70.790816326530
$ws.Range("N17").Value = -24.985994397759
$ws.Range("C18").Value = 232
$ws.Range("D18").Value = 254
$ws.Range("E18").Value = -8.661417322834
$ws.Range("F18").Value = 986
$ws.Range("G18").Value = 1053
$ws.Range("H18").Value = -6.362773029439
$ws.Range("I18").Value = 2684
$ws.Range("J18").Value = 3076
$ws.Range("K18").Value = -12.743823146944
$ws.Range("L18").Value = -15.517784073024
$ws.Range("M18").Value = -28.139223560910
$ws.Range("N18").Value = -86.876589086641
$ws.Range("C19").Value = 900
$ws.Range("D19").Value = 902
$ws.Range("E19").Value = -0.221729490022
$ws.Range("F19").Value = 3496
$ws.Range("G19").Value = 3527
$ws.Range("H19").Value = -0.878933938191
$ws.Range("I19").Value = 9845
$ws.Range("J19").Value = 9905
$ws.Range("K19").Value = -0.605754669358
$ws.Range("L19").Value = -3.157584103875
$ws.Range("M19").Value = 40.462262804965
$ws.Range("N19").Value = -37.072547139661
$ws.Range("C20").Value = 274
$ws.Range("D20").Value = 278
$ws.Range("E20").Value = -1.438848920863
$ws.Range("F20").Value = 903
$ws.Range("G20").Value = 1039
$ws.Range("H20").Value = -13.089509143407
$ws.Range("I20").Value = 2695
$ws.Range("J20").Value = 2980
$ws.Range("K20").Value = -9.563758389261
$ws.Range("L20").Value = -4.466501240694
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -88.451814714830
$ws.Range("C21").Value = 2297
$ws.Range("D21").Value = 2198
$ws.Range("E21").Value = 4.504094631483
$ws.Range("F21").Value = 8631
$ws.Range("G21").Value = 8637
$ws.Range("H21").Value = -0.069468565474
$ws.Range("I21").Value = 24296
$ws.Range("J21").Value = 24677
$ws.Range("K21").Value = -1.543947805648
$ws.Range("L21").Value = -0.625792465949
$ws.Range("M21").Value = 22.874627016638
$ws.Range("N21").Value = -71.466153051158
$ws.Range("C22").Value = 30
$ws.Range("D22").Value = 52
$ws.Range("E22").Value = -42.307692307692
$ws.Range("F22").Value = 148
$ws.Range("G22").Value = 169
$ws.Range("H22").Value = -12.426035502958
$ws.Range("I22").Value = 467
$ws.Range("J22").Value = 437
$ws.Range("K22").Value = 6.864988558352
$ws.Range("L22").Value = -7.157057654075
$ws.Range("M22").Value = 10.926365795724
$ws.Range("C23").Value = 118
$ws.Range("D23").Value = 109
$ws.Range("E23").Value = 8.256880733944
$ws.Range("F23").Value = 426
$ws.Range("H23").Value = -3.181818181818
$ws.Range("I23").Value = 1193
$ws.Range("J23").Value = 1232
$ws.Range("K23").Value = -3.165584415584
$ws.Range("L23").Value = 3.022452504317
$ws.Range("M23").Value = 58.643617021276
$ws.Range("C24").Value = 2014
$ws.Range("D24").Value = 1878
$ws.Range("E24").Value = 7.241746538871
$ws.Range("F24").Value = 8169
$ws.Range("G24").Value = 7829
$ws.Range("H24").Value = 4.342827947375
$ws.Range("I24").Value = 22250
$ws.Range("J24").Value = 21952
$ws.Range("K24").Value = 1.357507288629
$ws.Range("L24").Value = 4.942929912272
$ws.Range("M24").Value = 50.449658529988
$ws.Range("C25").Value = 1115
$ws.Range("D25").Value = 956
$ws.Range("E25").Value = 16.631799163179
$ws.Range("F25").Value = 4605
$ws.Range("G25").Value = 4057
$ws.Range("H25").Value = 13.507517870347
$ws.Range("I25").Value = 12450
$ws.Range("J25").Value = 11805
$ws.Range("K25").Value = 5.463786531130
$ws.Range("L25").Value = 3.508480212836
$ws.Range("C26").Value = 1015
$ws.Range("D26").Value = 785
$ws.Range("E26").Value = 29.299363057324
$ws.Range("F26").Value = 3539
$ws.Range("G26").Value = 2918
$ws.Range("H26").Value = 21.281699794379
$ws.Range("I26").Value = 9141
$ws.Range("J26").Value = 8165
$ws.Range("K26").Value = 11.953459889773
$ws.Range("L26").Value = 17.538896746817
$ws.Range("M26").Value = 6.364905748196
$ws.Range("C27").Value = 54
$ws.Range("D27").Value = 47
$ws.Range("E27").Value = 14.893617021276
$ws.Range("F27").Value = 191
$ws.Range("G27").Value = 158
$ws.Range("H27").Value = 20.886075949367
$ws.Range("I27").Value = 492
$ws.Range("J27").Value = 499
$ws.Range("K27").Value = -1.402805611222
$ws.Range("L27").Value = -12.142857142857
$ws.Range("C28").Value = 110
$ws.Range("D28").Value = 88
$ws.Range("E28").Value = 25
$ws.Range("F28").Value = 384
$ws.Range("G28").Value = 351
$ws.Range("H28").Value = 9.401709401709
$ws.Range("I28").Value = 986
$ws.Range("J28").Value = 977
$ws.Range("K28").Value = 0.921187308085
$ws.Range("L28").Value = 6.941431670282
$ws.Range("C29").Value = 23
$ws.Range("D29").Value = 24
$ws.Range("E29").Value = -4.166666666666
$ws.Range("F29").Value = 58
$ws.Range("G29").Value = 71
$ws.Range("H29").Value = -18.309859154929
$ws.Range("I29").Value = 180
$ws.Range("J29").Value = 219
$ws.Range("K29").Value = -17.808219178082
$ws.Range("L29").Value = -29.411764705882
$ws.Range("M29").Value = -31.297709923664
$ws.Range("N29").Value = -83.651226158038
$ws.Range("C30").Value = 19
$ws.Range("D30").Value = 20
$ws.Range("E30").Value = -5
$ws.Range("F30").Value = 52
$ws.Range("G30").Value = 59
$ws.Range("H30").Value = -11.864406779661
$ws.Range("I30").Value = 152
$ws.Range("J30").Value = 181
$ws.Range("K30").Value = -16.022099447513
$ws.Range("L30").Value = -33.624454148471
$ws.Range("M30").Value = -30.909090909090
$ws.Range("N30").Value = -84.965380811078
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 12
$ws.Range("E31").Value = -66.666666666666
$ws.Range("F31").Value = 43
$ws.Range("G31").Value = 37
$ws.Range("H31").Value = 16.216216216216
$ws.Range("I31").Value = 98
$ws.Range("J31").Value = 94
$ws.Range("K31").Value = 4.255319148936
$ws.Range("L31").Value = -42.690058479532
$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 3
$ws.Range("E33").Value = 133.333333333333
$ws.Range("F33").Value = 25
$ws.Range("H33").Value = 38.888888888888
$ws.Range("I33").Value = 56
$ws.Range("J33").Value = 43
$ws.Range("K33").Value = 30.232558139534
$ws.Range("L33").Value = 5.660377358490
